{"js": "// Apply the benchmark-table value refresh described by the commit:\n// \"Fixed README.md stats and docx preparation for all DaCapo - JDK 21 -\n// Shenandoah GC tests\".\n//\n// The document is a single table with one column and 46 rows, each row\n// holding one run of text in its only paragraph. The row count does not\n// change; only the text of specific rows is corrected. In particular,\n// the last three rows used to pack several tab-separated numbers into a\n// single run - the corrected run is a lone number (the new first-column\n// figure), matching what the first few rows now show.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-based row index -> corrected cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"157\",\n  4: \"0.00001\",\n  5: \"0.00049\",\n  6: \"0.00013\",\n  8: \"0.00022\",\n  9: \"0.00025\",\n  10: \"0.00028\",\n  11: \"0.02419\",\n  43: \"99.95\",\n  44: \"0.02\",\n  45: \"51\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table value refresh described by the commit:\n# \"Fixed README.md stats and docx preparation for all DaCapo - JDK 21 -\n# Shenandoah GC tests\". The document contains a single, single-column\n# table (46 rows x 1 column). Only the rows below actually change text;\n# the last three rows collapse their old multi-run, tab-separated\n# content down to a single number (their new first-column value), while\n# rows 1-3 and 7 in the table pick up fresh figures.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"157\"\n$t.Cell(5, 1).Range.Text = \"0.00001\"\n$t.Cell(6, 1).Range.Text = \"0.00049\"\n$t.Cell(7, 1).Range.Text = \"0.00013\"\n$t.Cell(9, 1).Range.Text = \"0.00022\"\n$t.Cell(10, 1).Range.Text = \"0.00025\"\n$t.Cell(11, 1).Range.Text = \"0.00028\"\n$t.Cell(12, 1).Range.Text = \"0.02419\"\n$t.Cell(44, 1).Range.Text = \"99.95\"\n$t.Cell(45, 1).Range.Text = \"0.02\"\n$t.Cell(46, 1).Range.Text = \"51\"\n"}
